# Apply "backup" column (R) addition and related data updates to the
# AXISBANK.NS 1mo stock-history sheet, plus append 6 new monthly rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cell R1 = "backup" (copy style from Q1 header, then set value)
# ---------------------------------------------------------------------
$ws.Range("Q1").Copy($ws.Range("R1"))
$ws.Range("R1").Value = "backup"

# ---------------------------------------------------------------------
# 2. Populate column R ("backup") for existing data rows 2-308.
#    R mirrors column P ("two_line_structure") for every row; for the
#    overwhelming majority P is 0, with a handful of exceptions.
# ---------------------------------------------------------------------
$rRow = 18   # column R index

for ($i = 2; $i -le 308; $i++) {
    $ws.Cells.Item($i, $rRow).Value = 0
}

# Rows where P (and therefore the new R/backup value) is non-zero
$nonZeroBackup = @{
    48  = 2
    55  = 2
    96  = 1
    105 = 1
    171 = 1
    183 = 1
}
foreach ($r in $nonZeroBackup.Keys) {
    $ws.Cells.Item($r, $rRow).Value = $nonZeroBackup[$r]
}

# ---------------------------------------------------------------------
# 3. Column Q ("detect_structure") recalculation: rows 13-44 (which sit
#    before the first P-checkpoint at row 48) had their stale detection
#    values reset to 0.
# ---------------------------------------------------------------------
$qResetRows = @(13, 17, 20, 25, 28, 31, 34, 39, 44)
foreach ($r in $qResetRows) {
    $ws.Cells.Item($r, 17).Value = 0
}

# ---------------------------------------------------------------------
# 4. Column O ("isPivot") recalculation affecting row 306 only, now that
#    additional rows extend the series.
# ---------------------------------------------------------------------
$ws.Cells.Item(306, 15).Value = 2

# ---------------------------------------------------------------------
# 5. Append six new monthly rows (309-314), continuing the data set
#    through December 2024. Column A keeps the date-number style; F
#    (Adj Close) and R (backup) are left blank for these newest rows.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=309; A=45474; B=1272.018267104002;  C=1338.616891265523; D=1153.110039464272; E=1165.20068359375;   G=244274229; H=2024; I=7;  N=27; O=1; Q=0 }
    @{ Row=310; A=45505; B=1168;                C=1184.849975585938; D=1123.099975585938; E=1175.25;            G=160196415; H=2024; I=8;  N=31; O=0; Q=1 }
    @{ Row=311; A=45536; B=1176;                C=1281.650024414062; D=1145;               E=1232.199951171875; G=166160797; H=2024; I=9;  N=35; O=0; Q=0 }
    @{ Row=312; A=45566; B=1228.099975585938;   C=1242.949951171875; D=1124;               E=1159.550048828125; G=220990545; H=2024; I=10; N=40; O=0; Q=0 }
    @{ Row=313; A=45597; B=1164.050048828125;   C=1187;               D=1115.75;            E=1136.300048828125; G=166280243; H=2024; I=11; N=44; O=0; Q=2 }
    @{ Row=314; A=45627; B=1132.699951171875;   C=1193.849975585938; D=1063.949951171875;  E=1069.949951171875; G=151148712; H=2024; I=12; N=48; O=0; Q=0 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row

    # Column A needs the same date-number format as the other rows.
    $ws.Range("A308").Copy($ws.Range("A$r"))
    $ws.Cells.Item($r, 1).Value = $nr.A

    $ws.Cells.Item($r, 2).Value = $nr.B    # Open
    $ws.Cells.Item($r, 3).Value = $nr.C    # High
    $ws.Cells.Item($r, 4).Value = $nr.D    # Low
    $ws.Cells.Item($r, 5).Value = $nr.E    # Close
    # F (Adj Close) intentionally left blank for these rows
    $ws.Cells.Item($r, 7).Value = $nr.G    # Volume
    $ws.Cells.Item($r, 8).Value = $nr.H    # Year
    $ws.Cells.Item($r, 9).Value = $nr.I    # Month
    $ws.Cells.Item($r, 10).Value = 1       # Day
    $ws.Cells.Item($r, 11).Value = 0       # Hour
    $ws.Cells.Item($r, 12).Value = 0       # Minute
    $ws.Cells.Item($r, 13).Value = 0       # Second
    $ws.Cells.Item($r, 14).Value = $nr.N   # Week
    $ws.Cells.Item($r, 15).Value = $nr.O   # isPivot
    $ws.Cells.Item($r, 16).Value = 0       # two_line_structure
    $ws.Cells.Item($r, 17).Value = $nr.Q   # detect_structure
    # R (backup) intentionally left blank for these rows
}
